$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect to allow cell edits, then restore protection at the end.
$ws.Unprotect()

# Update the confidential disclosure date (2021-05-11 -> 2021-05-12) in the A41 note cell.
$ws.Range("A41").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-05-12 for illustrative purposes only and are subject to change."

# Refresh Weight (D) and Percent Change (E) values for each holding row.
$ws.Range("D2").Value = 0.03165004432297032
$ws.Range("E2").Value = 0.0006092784402471629
$ws.Range("D3").Value = 0.02864159196696586
$ws.Range("E3").Value = -0.02699973206372763
$ws.Range("D4").Value = 0.02872423635113237
$ws.Range("E4").Value = -0.0285044904334244
$ws.Range("D5").Value = 0.06343763251387491
$ws.Range("E5").Value = -0.02232382417623313
$ws.Range("D6").Value = 0.01584528971284028
$ws.Range("E6").Value = -0.009164741822517142
$ws.Range("D7").Value = 0.01515579942207963
$ws.Range("E7").Value = -0.02921243281140451
$ws.Range("D8").Value = 0.02978739732173164
$ws.Range("E8").Value = -0.005826397146254547
$ws.Range("D9").Value = 0.03456502981592931
$ws.Range("E9").Value = -0.04019127860639871
$ws.Range("D10").Value = 0.02922541550939932
$ws.Range("E10").Value = -0.02040074331421182
$ws.Range("D11").Value = 0.03049302294463908
$ws.Range("E11").Value = -0.006685337428855243
$ws.Range("D12").Value = 0.01142165066415585
$ws.Range("E12").Value = -0.05357911964854856
$ws.Range("D13").Value = 0.01478783514019538
$ws.Range("E13").Value = -0.07749627421758554
$ws.Range("D14").Value = 0.01461565933984846
$ws.Range("E14").Value = -0.03581189331825452
$ws.Range("D15").Value = 0.009003515337912226
$ws.Range("E15").Value = -0.01774630649532305
$ws.Range("D16").Value = 0.008006469874646179
$ws.Range("E16").Value = -0.0238885202388851
$ws.Range("D17").Value = 0.02969137241822387
$ws.Range("E17").Value = -0.02328817962516239
$ws.Range("D18").Value = 0.02580944763051663
$ws.Range("E18").Value = -0.02046293190204618
$ws.Range("D19").Value = 0.03281277209925591
$ws.Range("E19").Value = 0.0011693802284789
$ws.Range("D20").Value = 0.03015831318876469
$ws.Range("E20").Value = -0.01298404723844315
$ws.Range("D21").Value = 0.0454300115210207
$ws.Range("E21").Value = -0.03018070306138376
$ws.Range("D22").Value = 0.03541115089191984
$ws.Range("E22").Value = -0.01533674149811082
$ws.Range("D23").Value = 0.03253237151011951
$ws.Range("E23").Value = -0.04131131676041866
$ws.Range("D24").Value = 0.03124666101930041
$ws.Range("E24").Value = -0.018476535932895
$ws.Range("D25").Value = 0.01439429045368815
$ws.Range("E25").Value = -0.03973917322834652
$ws.Range("D26").Value = 0.01475792574402083
$ws.Range("E26").Value = -0.03271999999999997
$ws.Range("D27").Value = 0.03121635807843936
$ws.Range("E27").Value = -0.0202468450977672
$ws.Range("D28").Value = 0.03115929409889581
$ws.Range("E28").Value = -0.004395271294331593
$ws.Range("D29").Value = 0.02907075244760198
$ws.Range("E29").Value = -0.02936279088656935
$ws.Range("D30").Value = 0.02967956607762866
$ws.Range("E30").Value = -0.02720245040840152
$ws.Range("D31").Value = 0.03378089202806367
$ws.Range("E31").Value = -0.03828746177370035
$ws.Range("D32").Value = 0.03152588097437729
$ws.Range("E32").Value = -0.009749399244764922
$ws.Range("D33").Value = 0.028867486617021
$ws.Range("E33").Value = -0.04111652636242791
$ws.Range("D34").Value = 0.03240446948700466
$ws.Range("E34").Value = -0.0153752732572261
$ws.Range("D35").Value = 0.03055599009448024
$ws.Range("E35").Value = -0.000115915150110002
$ws.Range("D36").Value = 0.03119628729942749
$ws.Range("E36").Value = -0.006875236533366857
$ws.Range("D37").Value = 0.03293811608190846
$ws.Range("E37").Value = -0.03001338176256929
$ws.Range("E38").Value = -0.02195782578367045

# Restore sheet protection (best-effort; original password hash cannot be re-derived here).
$ws.Protect()
